$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (week 3 / SEMANA=3) of data to add, rows 58-85, columns A-G
$data = @(
    @(2026, "ENERO", 3, "AMARILLO", "COLORES", "GOLDFINCH", 10310),
    @(2026, "ENERO", 3, "AMARILLO", "COLORES", "HIGH AND EXOTIC", 22195),
    @(2026, "ENERO", 3, "AMARILLO", "COLORES", "MOMENTUM", 9420),
    @(2026, "ENERO", 3, "AMARILLO", "COLORES", "SUNDAY MORNING", 1230),
    @(2026, "ENERO", 3, "BIC. AMARILLO", "COLORES", "SUMMER LIGHT", 4450),
    @(2026, "ENERO", 3, "BIC. ROSADO", "COLORES", "DISCOVERY", 500),
    @(2026, "ENERO", 3, "BICOLOR", "COLORES", "BLUSH", 3770),
    @(2026, "ENERO", 3, "BLANCO", "COLORES", "HIGH AND PURE", 14800),
    @(2026, "ENERO", 3, "BLANCO", "COLORES", "SUGGAR DOLL", 6300),
    @(2026, "ENERO", 3, "BLANCO", "COLORES", "VANILLA ICE", 3200),
    @(2026, "ENERO", 3, "BLANCO", "COLORES", "VENDELA", 7265),
    @(2026, "ENERO", 3, "DURAZNO", "COLORES", "TIFANY", 4082),
    @(2026, "ENERO", 3, "HOT PINK", "COLORES", "COTTON CANDY", 3120),
    @(2026, "ENERO", 3, "HOT PINK", "COLORES", "JACARANDA", 21350),
    @(2026, "ENERO", 3, "HOT PINK", "COLORES", "PINK FLOYD", 10496),
    @(2026, "ENERO", 3, "LAVANDER", "COLORES", "DEEP PURPLE", 22350),
    @(2026, "ENERO", 3, "LAVANDER", "COLORES", "MOODY BLUES", 13841),
    @(2026, "ENERO", 3, "NARANJA", "COLORES", "ALIVE", 6960),
    @(2026, "ENERO", 3, "NARANJA", "COLORES", "BROMO", 3950),
    @(2026, "ENERO", 3, "NARANJA", "COLORES", "CLEMENTINA", 2800),
    @(2026, "ENERO", 3, "NARANJA", "COLORES", "NINA", 23195),
    @(2026, "ENERO", 3, "ROJO", "ROJO", "FREEDOM", 228666),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "ABSOLUT IN PINK", 6830),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "HIGH AND BONITA", 16550),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "LUCIANO", 1880),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "PECKOUBO", 7977),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "STARFISH", 4417),
    @(2026, "ENERO", 3, "ROSADO", "COLORES", "TABATHA", 14359)
)

$startRow = 58
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}

# Update the view to match: scrolled down, with the active selection on A85
$excel.ActiveWindow.ScrollRow = 53
$ws.Range("A85").Select() | Out-Null
